$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30, shifting existing rows 30-54 down to 31-55
$ws.Rows("30:30").Insert()

# Populate the new row 30 with the new weekly record
$ws.Range("A30").Value = 8
$ws.Range("B30").Value = "Terminal La Palmera de La Serena"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").Value = 44512
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = 100112052
$ws.Range("G30").Value = "Albahaca"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 800
$ws.Range("K30").Value = 3000
$ws.Range("L30").Value = 3500
$ws.Range("M30").Value = 3250
$ws.Range("N30").Value = "`$/paquete"
$ws.Range("O30").Value = "Región de Arica y Parinacota"
$ws.Range("P30").Value = 3250
$ws.Range("Q30").Value = 1
$ws.Range("R30").Value = "Hortaliza"
